$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2433.5
$ws.Range("I62").Value = 2444
$ws.Range("J62").Value = 2402
$ws.Range("K62").Value = 2444
$ws.Range("L62").Value = 2402
$ws.Range("M62").Value = -1820
$ws.Range("N62").Value = -3650
$ws.Range("H65").Value = 2433.5
$ws.Range("I65").Value = 2444
$ws.Range("J65").Value = 2402
$ws.Range("K65").Value = 12220
$ws.Range("L65").Value = 12010
$ws.Range("M65").Value = -9100
$ws.Range("N65").Value = -18250
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H125").Value = 1574.2667
$ws.Range("I125").Value = 1201.75
$ws.Range("J125").Value = 2000
$ws.Range("K125").Value = 10815.75
$ws.Range("L125").Value = 18000
$ws.Range("M125").Value = -8355.75
$ws.Range("N125").Value = -22920
$ws.Range("H132").Value = 4610095
$ws.Range("I132").Value = 5496328.5
$ws.Range("K132").Value = 16488985.5
$ws.Range("M132").Value = -16486455.5
$ws.Range("H137").Value = 2299.3333
$ws.Range("I137").Value = 2333.1667
$ws.Range("K137").Value = 6999.500100000001
$ws.Range("M137").Value = -4449.500100000001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22833.32
$ws.Range("I32").Value = 24457.652
$ws.Range("K32").Value = 24457.652
$ws.Range("M32").Value = -24170.652
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 4719.1035
$ws.Range("I132").Value = 4674.16
$ws.Range("K132").Value = 14022.48
$ws.Range("M132").Value = -11492.48

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1949
$ws.Range("I86").Value = 1939.7273
$ws.Range("K86").Value = 1939.7273
$ws.Range("M86").Value = -816.7273
$ws.Range("H89").Value = 1949
$ws.Range("I89").Value = 1939.7273
$ws.Range("K89").Value = 9698.636500000001
$ws.Range("M89").Value = -4082.636500000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 15580
$ws.Range("J36").Value = 15996
$ws.Range("L36").Value = 15996
$ws.Range("N36").Value = -16772
$ws.Range("H40").Value = 15580
$ws.Range("J40").Value = 15996
$ws.Range("L40").Value = 15996
$ws.Range("N40").Value = -16316
$ws.Range("H58").Value = 1172
$ws.Range("I58").Value = 1244.7693
$ws.Range("J58").Value = 982.8
$ws.Range("K58").Value = 1244.7693
$ws.Range("L58").Value = 982.8
$ws.Range("M58").Value = -1041.7693
$ws.Range("N58").Value = -1388.8
$ws.Range("H99").Value = 1542.5143
$ws.Range("I99").Value = 1135.8182
$ws.Range("J99").Value = 2230.7693
$ws.Range("K99").Value = 1135.8182
$ws.Range("L99").Value = 2230.7693
$ws.Range("M99").Value = 362.1818000000001
$ws.Range("N99").Value = -5226.7693
$ws.Range("H107").Value = 506.48148
$ws.Range("I107").Value = 524.05884
$ws.Range("J107").Value = 476.6
$ws.Range("K107").Value = 524.05884
$ws.Range("L107").Value = 476.6
$ws.Range("M107").Value = 1395.94116
$ws.Range("N107").Value = -4316.6
$ws.Range("H126").Value = 1542.5143
$ws.Range("I126").Value = 1135.8182
$ws.Range("J126").Value = 2230.7693
$ws.Range("K126").Value = 3407.4546
$ws.Range("L126").Value = 6692.3079
$ws.Range("M126").Value = -937.4546
$ws.Range("N126").Value = -11632.3079
$ws.Range("H136").Value = 1172
$ws.Range("I136").Value = 1244.7693
$ws.Range("J136").Value = 982.8
$ws.Range("K136").Value = 3734.3079
$ws.Range("L136").Value = 2948.4
$ws.Range("M136").Value = -1184.3079
$ws.Range("N136").Value = -8048.4

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1560.4706
$ws.Range("I34").Value = 689.5
$ws.Range("J34").Value = 1676.6
$ws.Range("K34").Value = 2068.5
$ws.Range("L34").Value = 5029.799999999999
$ws.Range("M34").Value = -1984.5
$ws.Range("N34").Value = -5197.799999999999
$ws.Range("H39").Value = 2083.889
$ws.Range("J39").Value = 2083.889
$ws.Range("L39").Value = 6251.667
$ws.Range("N39").Value = -6839.667
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 3259.6
$ws.Range("I5").Value = 3259.6
$ws.Range("K5").Value = 3259.6
$ws.Range("M5").Value = -3147.6
$ws.Range("H70").Value = 8099814.5
$ws.Range("I70").Value = 11091152
$ws.Range("J70").Value = 5605.294
$ws.Range("K70").Value = 11091152
$ws.Range("L70").Value = 5605.294
$ws.Range("M70").Value = -11090882
$ws.Range("N70").Value = -6145.294
$ws.Range("H73").Value = 8099814.5
$ws.Range("I73").Value = 11091152
$ws.Range("J73").Value = 5605.294
$ws.Range("K73").Value = 11091152
$ws.Range("L73").Value = 5605.294
$ws.Range("M73").Value = -11090216
$ws.Range("N73").Value = -7477.294
$ws.Range("H122").Value = 1912.2354
$ws.Range("I122").Value = 666.6667
$ws.Range("J122").Value = 2591.6365
$ws.Range("K122").Value = 2000.0001
$ws.Range("L122").Value = 7774.9095
$ws.Range("M122").Value = 449.9999
$ws.Range("N122").Value = -12674.9095

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 5000
$ws.Range("J38").Value = 5000
$ws.Range("L38").Value = 5000
$ws.Range("N38").Value = -5820
$ws.Range("H40").Value = 2793.6667
$ws.Range("I40").Value = 2302.4
$ws.Range("K40").Value = 2302.4
$ws.Range("M40").Value = -2166.4
$ws.Range("H132").Value = 24243.783
$ws.Range("I132").Value = 32631.812
$ws.Range("J132").Value = 5071.143
$ws.Range("K132").Value = 97895.436
$ws.Range("L132").Value = 15213.429
$ws.Range("M132").Value = -95365.436
$ws.Range("N132").Value = -20273.429

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 849.73334
$ws.Range("I122").Value = 784.35
$ws.Range("J122").Value = 980.5
$ws.Range("K122").Value = 2353.05
$ws.Range("L122").Value = 2941.5
$ws.Range("M122").Value = 96.94999999999982
$ws.Range("N122").Value = -7841.5

